$d = $word.ActiveDocument

# 1. Narrow the first (label) column of the single table from 11638 -> 5376 dxa.
#    Word COM widths are expressed in points (1 point = 20 dxa/twips).
$t = $d.Tables.Item(1)
$t.Columns.Item(1).Width = 5376 / 20.0

# 2. Tighten the row heights of the two question-label rows from 619 -> 617 dxa.
$t.Rows.Item(6).Height = 617 / 20.0
$t.Rows.Item(10).Height = 617 / 20.0

# 3. Shorten/rename the two long Qualtrics question labels to their new short form.
#    Set the text directly on the label cell's Range so xml:space="preserve" /
#    run formatting (bold, font, color) on the single run is kept intact,
#    matching how Find/Replace would behave inside Word itself.
$t.Cell(6, 1).Range.Text = "Percent of recent studies that could be replicated"
$t.Cell(10, 1).Range.Text = "Percent of recent studies that should be replicated"
